$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Golang Developer", "https://www.dice.com/job-detail/420adaca-6f15-4ba9-83c3-f32f835793c4", "Hybrid in McLean, Virginia", "Contract, Third Party", "Depends on Experience", "NimbusAITech LLC"),
    @("Golang Developer (Backend)", "https://www.dice.com/job-detail/648c759a-755a-4cb4-928c-da48f86f2562", "Hybrid in Richmond, Virginia", "Contract, Third Party", "$60 - $70", "Spar Information Systems"),
    @("Golang developer with AWS (Hybrid in Richmond, VA)", "https://www.dice.com/job-detail/1e88c296-9335-4e25-9ee2-772117afa907", "Hybrid in Richmond, Virginia", "Contract, Third Party", "Depends on Experience", "InfiCare Technologies"),
    @("Golang + MQTT working", "https://www.dice.com/job-detail/ae99ab4a-c621-482b-b39f-4fb4f33353cd", "Atlanta, Georgia", "Contract", "Depends on Experience", "SUS Infotech Inc"),
    @("Onsite Role: Sr Java with Go & Devops exp (TX, 8+Months)", "https://www.dice.com/job-detail/e41efe3a-2480-4a0a-8fe2-8048ad2fedbc", "Texas", "Full-time, Part-time, Third Party, Contract", "`$DOE", "Elista Global LLC"),
    @("Golang Developer", "https://www.dice.com/job-detail/6425197d-aaa4-433a-9c76-4b28756c3846", "Plano, Texas", "Contract", "Depends on Experience", "iTech US, Inc."),
    @("Senior Golang Developer", "https://www.dice.com/job-detail/dd9ef790-30ab-4641-b94c-f6db77c26248", "Alpharetta, Georgia", "Contract", "Depends on Experience", "Swanktek")
)

$startRow = 92
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($row, $c).Value = $rec[$c - 1]
    }
}
